# Generate Report for handback
#
# Both localization targets (zh-cn and de-de) for
# 3aaff79e-7311-419c-9ef3-0ea864b799da.md have now been handed back (they
# were previously only "Ready for handoff"). Regenerating the status report
# swaps the row ordering (the 3aaff79e file now sorts/report before the
# 9eb1fb6a file) and refreshes the "Handed back" status + handback
# timestamps for that file.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de
# Row 2 and Row 3 swap which file they describe, and the status for the
# 3aaff79e file becomes "Handed back: in sync with en-US" in both columns.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("A2").Value() = "3aaff79e-7311-419c-9ef3-0ea864b799da.md"
$overview.Range("B2").Value() = "Handed back: in sync with en-US"
$overview.Range("C2").Value() = "Handed back: in sync with en-US"

$overview.Range("A3").Value() = "9eb1fb6a-3318-4339-a92b-8f71d363a4eb.md"
$overview.Range("B3").Value() = "Handed back: in sync with en-US"
$overview.Range("C3").Value() = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Sheet "zh-cn" detail table:
# Source File Name | Status | Latest Handoff File | Latest Handoff Datetime |
# Latest Target File | Latest Handback File | Latest Handback DateTime |
# Handoff Reason | Dependency From
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("A2").Value() = "3aaff79e-7311-419c-9ef3-0ea864b799da.md"
$zhcn.Range("B2").Value() = "Handed back: in sync with en-US"
$zhcn.Range("C2").Value() = "3aaff79e-7311-419c-9ef3-0ea864b799da.db4c6cffa82d6430ba90646cbcf924abcdd63c90.zh-cn.xlf"
$zhcn.Range("D2").Value() = "2016-02-19 05:46:43"
$zhcn.Range("E2").Value() = "3aaff79e-7311-419c-9ef3-0ea864b799da.md"
$zhcn.Range("F2").Value() = "3aaff79e-7311-419c-9ef3-0ea864b799da.db4c6cffa82d6430ba90646cbcf924abcdd63c90.zh-cn.xlf"
$zhcn.Range("G2").Value() = "2016-02-19 05:47:25"
$zhcn.Range("H2").Value() = "Include"

$zhcn.Range("A3").Value() = "9eb1fb6a-3318-4339-a92b-8f71d363a4eb.md"
$zhcn.Range("B3").Value() = "Handed back: in sync with en-US"
$zhcn.Range("C3").Value() = "9eb1fb6a-3318-4339-a92b-8f71d363a4eb.be3be81da801bc3ff874f4e44f79467f38bb3f5e.zh-cn.xlf"
$zhcn.Range("D3").Value() = "2016-02-19 05:44:46"
$zhcn.Range("E3").Value() = "9eb1fb6a-3318-4339-a92b-8f71d363a4eb.md"
$zhcn.Range("F3").Value() = "9eb1fb6a-3318-4339-a92b-8f71d363a4eb.be3be81da801bc3ff874f4e44f79467f38bb3f5e.zh-cn.xlf"
$zhcn.Range("G3").Value() = "2016-02-19 05:45:27"
$zhcn.Range("H3").Value() = "Include"

# ---------------------------------------------------------------------
# Sheet "de-de" detail table (same column layout as zh-cn).
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("A2").Value() = "3aaff79e-7311-419c-9ef3-0ea864b799da.md"
$dede.Range("B2").Value() = "Handed back: in sync with en-US"
$dede.Range("C2").Value() = "3aaff79e-7311-419c-9ef3-0ea864b799da.db4c6cffa82d6430ba90646cbcf924abcdd63c90.de-de.xlf"
$dede.Range("D2").Value() = "2016-02-19 05:46:53"
$dede.Range("E2").Value() = "3aaff79e-7311-419c-9ef3-0ea864b799da.md"
$dede.Range("F2").Value() = "3aaff79e-7311-419c-9ef3-0ea864b799da.db4c6cffa82d6430ba90646cbcf924abcdd63c90.de-de.xlf"
$dede.Range("G2").Value() = "2016-02-19 05:47:42"
$dede.Range("H2").Value() = "Include"

$dede.Range("A3").Value() = "9eb1fb6a-3318-4339-a92b-8f71d363a4eb.md"
$dede.Range("B3").Value() = "Handed back: in sync with en-US"
$dede.Range("C3").Value() = "9eb1fb6a-3318-4339-a92b-8f71d363a4eb.be3be81da801bc3ff874f4e44f79467f38bb3f5e.de-de.xlf"
$dede.Range("D3").Value() = "2016-02-19 05:44:57"
$dede.Range("E3").Value() = "9eb1fb6a-3318-4339-a92b-8f71d363a4eb.md"
$dede.Range("F3").Value() = "9eb1fb6a-3318-4339-a92b-8f71d363a4eb.be3be81da801bc3ff874f4e44f79467f38bb3f5e.de-de.xlf"
$dede.Range("G3").Value() = "2016-02-19 05:45:43"
$dede.Range("H3").Value() = "Include"
